$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample placeholder values in the "Total value" / "Total number"
# cells of both Form 3A blocks (rows 6 and 9).
$ws.Cells.Item(6, 2).ClearContents() | Out-Null
$ws.Cells.Item(6, 3).ClearContents() | Out-Null
$ws.Cells.Item(9, 2).ClearContents() | Out-Null
$ws.Cells.Item(9, 3).ClearContents() | Out-Null

# Give the "Total value of payment transactions" cells a currency number
# format.
$ws.Cells.Item(6, 2).NumberFormat = """$""#,##0.00"
$ws.Cells.Item(9, 2).NumberFormat = """$""#,##0.00"

# Match the author's saved selection / window position.
$ws.Range("G6").Select() | Out-Null

$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
